# Update the dSF (column F) values for a handful of rows based on a
# repull of the underlying data / recalculated means.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -14
$ws.Range("F4").Value = -5
$ws.Range("F5").Value = -2
$ws.Range("F8").Value = 1
